# Analysis as of Sept 2021 Submission
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demographics")

# --- RACE/ETHNICITY block: rows 12-14 get relabeled/reordered ---
# Row 12: native american -> pacific islander (counts unchanged)
$ws.Range("A12").Value = "pacific islander"

# Row 13: asian -> native american; D13 and F13 counts updated
$ws.Range("A13").Value = "native american"
$ws.Range("D13").Value = "1/62 (1.6%)"
$ws.Range("F13").Value = "1/171 (0.6%)"

# Row 14: pacific islander -> asian; D14 and F14 counts updated
$ws.Range("A14").Value = "asian"
$ws.Range("D14").Value = "0/62 (0.0%)"
$ws.Range("F14").Value = "2/171 (1.2%)"

# --- Severity block: rows 35-36 updated counts ---
# Row 35: mild
$ws.Range("B35").Value = "93/510 (18.2%)"
$ws.Range("E35").Value = "55/266 (20.7%)"
$ws.Range("F35").Value = "30/171 (17.5%)"

# Row 36: none
$ws.Range("B36").Value = "0/510 (0.0%)"
$ws.Range("E36").Value = "0/266 (0.0%)"
$ws.Range("F36").Value = "0/171 (0.0%)"
